# Edit script: insert two new rows (94-95) into the Perejil price sheet to add
# a new "semanal" observation block, shifting all subsequent rows down by 2.
#
# Net effect (per the supplied diff):
#   - Two blank rows are inserted before the current row 94, which pushes the
#     existing rows 94..177 down to 96..179 (preserving all of their values).
#   - The two newly-inserted rows (94 and 95) are populated with a new
#     Primera/Segunda price pair for Fecha serial 44880, using the same
#     Mercado/Region/etc. values as the template rows, with updated prices.
#   - The sheet dimension grows from A1:R177 to A1:R179 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, empty rows at position 94 (this shifts old rows 94-177 down
# to 96-179, exactly like the diff shows).
$ws.Range("A94:A95").EntireRow.Insert()

# --- New row 94 ("Primera") ---
$ws.Cells.Item(94, 1).Value = 11
$ws.Cells.Item(94, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(94, 3).Value = "Bíobío"
$ws.Cells.Item(94, 4).Value2 = 44880
$ws.Cells.Item(94, 5).Value = 8
$ws.Cells.Item(94, 6).Value = 100112044
$ws.Cells.Item(94, 7).Value = "Perejil"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 200
$ws.Cells.Item(94, 11).Value = 700
$ws.Cells.Item(94, 12).Value = 800
$ws.Cells.Item(94, 13).Value = 750
$ws.Cells.Item(94, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(94, 15).Value = "Región de Ñuble"
$ws.Cells.Item(94, 16).Value = 750
$ws.Cells.Item(94, 17).Value = 1
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# --- New row 95 ("Segunda") ---
$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).Value2 = 44880
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = 100112044
$ws.Cells.Item(95, 7).Value = "Perejil"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Segunda"
$ws.Cells.Item(95, 10).Value = 100
$ws.Cells.Item(95, 11).Value = 600
$ws.Cells.Item(95, 12).Value = 600
$ws.Cells.Item(95, 13).Value = 600
$ws.Cells.Item(95, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(95, 15).Value = "Región de Ñuble"
$ws.Cells.Item(95, 16).Value = 600
$ws.Cells.Item(95, 17).Value = 1
$ws.Cells.Item(95, 18).Value = "Hortaliza"

# Make sure the date column keeps the expected date number format (it is
# normally inherited automatically from the row above on Insert, but set it
# explicitly too, to be safe).
$ws.Range("D94:D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
